# Refresh the cryptos worksheet Price (D) / Volume(1h) (E) columns
# for rows 2-51, matching the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.652.34"
$ws.Range("E2").Value = "  +6.43%  "

$ws.Range("D3").Value = "2.745.41"
$ws.Range("E3").Value = "  +5.39%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "116.45"
$ws.Range("E5").Value = "  +6.43%  "

$ws.Range("D6").Value = "333.12"
$ws.Range("E6").Value = "  +3.74%  "

$ws.Range("D7").Value = "0.534"
$ws.Range("E7").Value = "  +2.53%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +5.90%  "

$ws.Range("D10").Value = "41.43"
$ws.Range("E10").Value = "  +5.72%  "

$ws.Range("D11").Value = "0.0861"
$ws.Range("E11").Value = "  +6.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.10"
$ws.Range("E12").Value = "  +1.94%  "

$ws.Range("E13").Value = "  +2.88%  "

$ws.Range("E14").Value = "  +4.95%  "

$ws.Range("D15").Value = "3.175.24"
$ws.Range("E15").Value = "  +5.52%  "

$ws.Range("D16").Value = "2.741.37"
$ws.Range("E16").Value = "  +5.46%  "

$ws.Range("D17").Value = "0.879"
$ws.Range("E17").Value = "  +2.22%  "

$ws.Range("D18").Value = "51.582.73"
$ws.Range("E18").Value = "  +6.46%  "

$ws.Range("D19").Value = "3.13"
$ws.Range("E19").Value = "  +6.50%  "

$ws.Range("D20").Value = "13.48"
$ws.Range("E20").Value = "  +5.43%  "

$ws.Range("D21").Value = "6.82"
$ws.Range("E21").Value = "  +2.56%  "

$ws.Range("D22").Value = "0.0₃0978"
$ws.Range("E22").Value = "  +3.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "279.50"
$ws.Range("E23").Value = "  +3.93%  "

$ws.Range("D24").Value = "69.46"
$ws.Range("E24").Value = "  +1.44%  "

$ws.Range("D25").Value = "2.65"
$ws.Range("E25").Value = "  +5.17%  "

$ws.Range("D26").Value = "26.69"
$ws.Range("E26").Value = "  +2.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").Value = "10.17"
$ws.Range("E28").Value = "  +1.92%  "

$ws.Range("E29").Value = "  +0.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.140"
$ws.Range("E30").Value = "  +2.17%  "

$ws.Range("D31").Value = "34.99"
$ws.Range("E31").Value = "  +0.73%  "

$ws.Range("D32").Value = "49.98"
$ws.Range("E32").Value = "  +1.57%  "

$ws.Range("D33").Value = "5.54"
$ws.Range("E33").Value = "  +1.75%  "

$ws.Range("D34").Value = "0.0817"
$ws.Range("E34").Value = "  +2.90%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").Value = "18.94"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").Value = "4.98"
$ws.Range("E37").Value = "  +1.01%  "

$ws.Range("E38").Value = "  +2.32%  "

$ws.Range("D39").Value = "3.16"
$ws.Range("E39").Value = "  +1.19%  "

$ws.Range("D40").Value = "128.11"
$ws.Range("E40").Value = "  +1.72%  "

$ws.Range("E41").Value = "  +9.50%  "

$ws.Range("D42").Value = "22.96"
$ws.Range("E42").Value = "  +4.24%  "

$ws.Range("E43").Value = "  +2.60%  "

$ws.Range("D44").Value = "2.26"
$ws.Range("E44").Value = "  +6.73%  "

$ws.Range("E45").Value = "  +12.83%  "

$ws.Range("D46").Value = "2.092.67"
$ws.Range("E46").Value = "  +1.83%  "

$ws.Range("D47").Value = "3.31"
$ws.Range("E47").Value = "  +3.04%  "

$ws.Range("D48").Value = "2.23"
$ws.Range("E48").Value = "  +2.91%  "

$ws.Range("D49").Value = "5.53"
$ws.Range("E49").Value = "  +7.82%  "

$ws.Range("D50").Value = "8.93"
$ws.Range("E50").Value = "  +0.85%  "

$ws.Range("D51").Value = "59.77"
$ws.Range("E51").Value = "  +2.37%  "
